$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H34").Value = 3650
$ws.Range("I34").Value = 3650
$ws.Range("K34").Value = 3650
$ws.Range("M34").Value = -3447
$ws.Range("H36").Value = 3650
$ws.Range("I36").Value = 3650
$ws.Range("K36").Value = 3650
$ws.Range("M36").Value = -2935
$ws.Range("H53").Value = 221.32259
$ws.Range("I53").Value = 179.57143
$ws.Range("J53").Value = 255.70589
$ws.Range("K53").Value = 179.57143
$ws.Range("L53").Value = 255.70589
$ws.Range("M53").Value = 457.42857
$ws.Range("N53").Value = -1529.70589
$ws.Range("H113").Value = 9302.6875
$ws.Range("I113").Value = 5879
$ws.Range("K113").Value = 5879
$ws.Range("M113").Value = -2625
$ws.Range("H132").Value = 13098.911
$ws.Range("I132").Value = 1670.0488
$ws.Range("K132").Value = 5010.1464
$ws.Range("M132").Value = -2480.1464
$ws.Range("H135").Value = 12827736
$ws.Range("I135").Value = 19235036
$ws.Range("K135").Value = 173115324
$ws.Range("M135").Value = -173112789
$ws.Range("H138").Value = 3927.65
$ws.Range("I138").Value = 2631.4443
$ws.Range("J138").Value = 4988.1816
$ws.Range("K138").Value = 7894.3329
$ws.Range("L138").Value = 14964.5448
$ws.Range("M138").Value = -2754.3329
$ws.Range("N138").Value = -25244.5448

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5057.5513
$ws.Range("I32").Value = 5308.233
$ws.Range("J32").Value = 1397.6
$ws.Range("K32").Value = 5308.233
$ws.Range("L32").Value = 1397.6
$ws.Range("M32").Value = -5021.233
$ws.Range("N32").Value = -1971.6
$ws.Range("H45").Value = 4321.8438
$ws.Range("I45").Value = 3456.28
$ws.Range("J45").Value = 7413.143
$ws.Range("K45").Value = 3456.28
$ws.Range("L45").Value = 7413.143
$ws.Range("M45").Value = -3079.28
$ws.Range("N45").Value = -8167.143
$ws.Range("H61").Value = 950
$ws.Range("I61").Value = 900
$ws.Range("J61").Value = 1000
$ws.Range("K61").Value = 900
$ws.Range("L61").Value = 1000
$ws.Range("M61").Value = -688
$ws.Range("N61").Value = -1424
$ws.Range("H68").Value = 45000
$ws.Range("J68").Value = 45000
$ws.Range("L68").Value = 45000
$ws.Range("N68").Value = -46622
$ws.Range("H71").Value = 45000
$ws.Range("J71").Value = 45000
$ws.Range("L71").Value = 135000
$ws.Range("N71").Value = -143112
$ws.Range("H110").Value = 4063.6584
$ws.Range("I110").Value = 4364.4
$ws.Range("K110").Value = 4364.4
$ws.Range("M110").Value = -2319.4
$ws.Range("H136").Value = 950
$ws.Range("I136").Value = 900
$ws.Range("J136").Value = 1000
$ws.Range("K136").Value = 2700
$ws.Range("L136").Value = 3000
$ws.Range("M136").Value = -150
$ws.Range("N136").Value = -8100
$ws.Range("H139").Value = 89714
$ws.Range("J139").Value = 89714
$ws.Range("L139").Value = 89714
$ws.Range("N139").Value = -99994

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H81").Value = 54650.8
$ws.Range("J81").Value = 58422.11
$ws.Range("L81").Value = 58422.11
$ws.Range("N81").Value = -60544.11
$ws.Range("H84").Value = 54650.8
$ws.Range("J84").Value = 58422.11
$ws.Range("L84").Value = 175266.33
$ws.Range("N84").Value = -185874.33
$ws.Range("H86").Value = 2226
$ws.Range("I86").Value = 1833
$ws.Range("J86").Value = 2619
$ws.Range("K86").Value = 1833
$ws.Range("L86").Value = 2619
$ws.Range("M86").Value = -710
$ws.Range("N86").Value = -4865
$ws.Range("H89").Value = 2226
$ws.Range("I89").Value = 1833
$ws.Range("J89").Value = 2619
$ws.Range("K89").Value = 9165
$ws.Range("L89").Value = 13095
$ws.Range("M89").Value = -3549
$ws.Range("N89").Value = -24327
$ws.Range("H134").Value = 7982.4
$ws.Range("I134").Value = 2602.8076
$ws.Range("J134").Value = 42949.75
$ws.Range("K134").Value = 7808.4228
$ws.Range("L134").Value = 128849.25
$ws.Range("M134").Value = -5273.4228
$ws.Range("N134").Value = -133919.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 8354.647000000001
$ws.Range("I31").Value = 33954
$ws.Range("K31").Value = 33954
$ws.Range("M31").Value = -33659
$ws.Range("H34").Value = 8354.647000000001
$ws.Range("I34").Value = 33954
$ws.Range("K34").Value = 33954
$ws.Range("M34").Value = -33752
$ws.Range("H140").Value = 81285.71000000001
$ws.Range("J140").Value = 81285.71000000001
$ws.Range("L140").Value = 81285.71000000001
$ws.Range("N140").Value = -91645.71000000001
$ws.Range("H141").Value = 36817.668
$ws.Range("J141").Value = 36817.668
$ws.Range("L141").Value = 36817.668
$ws.Range("N141").Value = -47177.668

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 459.8
$ws.Range("I11").Value = 549.75
$ws.Range("J11").Value = 100
$ws.Range("K11").Value = 1649.25
$ws.Range("L11").Value = 300
$ws.Range("M11").Value = -1509.25
$ws.Range("N11").Value = -580
$ws.Range("H136").Value = 3806.25
$ws.Range("J136").Value = 4166.6665
$ws.Range("L136").Value = 12499.9995
$ws.Range("N136").Value = -22699.9995

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 2005
$ws.Range("I5").Value = 0
$ws.Range("K5").Value = 0
$ws.Range("M5").Value = ""
$ws.Range("H70").Value = 228284
$ws.Range("I70").Value = 283356.75
$ws.Range("K70").Value = 283356.75
$ws.Range("M70").Value = -283086.75
$ws.Range("H73").Value = 228284
$ws.Range("I73").Value = 283356.75
$ws.Range("K73").Value = 283356.75
$ws.Range("M73").Value = -282420.75
$ws.Range("H122").Value = 5877.4736
$ws.Range("I122").Value = 5463.8887
$ws.Range("K122").Value = 16391.6661
$ws.Range("M122").Value = -13941.6661
$ws.Range("H126").Value = 4425.6
$ws.Range("I126").Value = 3948.5
$ws.Range("J126").Value = 4743.6665
$ws.Range("K126").Value = 11845.5
$ws.Range("L126").Value = 14230.9995
$ws.Range("M126").Value = -9375.5
$ws.Range("N126").Value = -19170.9995

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 5863.8887
$ws.Range("I40").Value = 6096.75
$ws.Range("J40").Value = 4001
$ws.Range("K40").Value = 6096.75
$ws.Range("L40").Value = 4001
$ws.Range("M40").Value = -5960.75
$ws.Range("N40").Value = -4273
$ws.Range("H55").Value = 332.5909
$ws.Range("I55").Value = 452.57144
$ws.Range("J55").Value = 122.625
$ws.Range("K55").Value = 452.57144
$ws.Range("L55").Value = 122.625
$ws.Range("M55").Value = -279.57144
$ws.Range("N55").Value = -468.625
$ws.Range("H61").Value = 1252
$ws.Range("J61").Value = 0
$ws.Range("L61").Value = 0
$ws.Range("N61").Value = ""
$ws.Range("H113").Value = 1252
$ws.Range("J113").Value = 0
$ws.Range("L113").Value = 0
$ws.Range("N113").Value = ""
$ws.Range("H136").Value = 4558.7617
$ws.Range("I136").Value = 6998.4
$ws.Range("J136").Value = 2340.9092
$ws.Range("K136").Value = 20995.2
$ws.Range("L136").Value = 7022.7276
$ws.Range("M136").Value = -18445.2
$ws.Range("N136").Value = -12122.7276

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H39").Value = 30000
$ws.Range("I39").Value = 0
$ws.Range("J39").Value = 30000
$ws.Range("K39").Value = 0
$ws.Range("L39").Value = 30000
$ws.Range("N39").Value = -30826
$ws.Range("M39").Value = ""
$ws.Range("H40").Value = 25165
$ws.Range("J40").Value = 30247.5
$ws.Range("L40").Value = 30247.5
$ws.Range("N40").Value = -30545.5
$ws.Range("H69").Value = 29826.445
$ws.Range("J69").Value = 29826.445
$ws.Range("L69").Value = 29826.445
$ws.Range("N69").Value = -31324.445
$ws.Range("H72").Value = 29826.445
$ws.Range("J72").Value = 29826.445
$ws.Range("L72").Value = 89479.33499999999
$ws.Range("N72").Value = -96967.33499999999
$ws.Range("H132").Value = 4831.973
$ws.Range("I132").Value = 4369.52
$ws.Range("J132").Value = 5795.4165
$ws.Range("K132").Value = 13108.56
$ws.Range("L132").Value = 17386.2495
$ws.Range("M132").Value = -10578.56
$ws.Range("N132").Value = -22446.2495
$ws.Range("H136").Value = 4754.3447
$ws.Range("I136").Value = 3697.4375
$ws.Range("J136").Value = 6055.154
$ws.Range("K136").Value = 11092.3125
$ws.Range("L136").Value = 18165.462
$ws.Range("M136").Value = -8542.3125
$ws.Range("N136").Value = -23265.462
